$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Header row: B1 used to hold "structure_string" (the raw PDB/SDF/
#    MOL2 text header) -- it is renamed to "structure_data" now that
#    the url/text columns have been merged into a single column.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "structure_data"

# ------------------------------------------------------------------
# 2. Rows whose "data" used to live in column C (the url rows) now
#    need that value moved into column B. Column B for these rows was
#    blank before, styled with the generic vertical-center style, so
#    reset it back to the default "Normal" style first.
# ------------------------------------------------------------------
$ws.Range("B3").Style = "Normal"
$ws.Range("B3").Value = "https://files.rcsb.org/view/5ZMZ.cif"

$ws.Range("B4").Style = "Normal"
$ws.Range("B4").Value = "https://files.rcsb.org/view/7P9W.pdb"

$ws.Range("B10").Value = "https://files.rcsb.org/view/2MCJ.cif"
$ws.Range("B10").Style = "Normal"

# ------------------------------------------------------------------
# 3. Row 2 (5jxe) keeps its hyperlink, but the hyperlink now lives on
#    column B instead of column C.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Range("B2").Value = "https://files.rcsb.org/view/5JXE.cif"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://files.rcsb.org/view/5JXE.cif")
$ws.Range("B2").Style = "Hyperlink"

# ------------------------------------------------------------------
# 4. Drop the old "structure_path" (C) and "structure_type" (D)
#    columns now that their contents have been folded into column B
#    (path) or removed entirely (type). The old "description" column
#    (E) slides left into column C automatically.
# ------------------------------------------------------------------
$ws.Columns("C:D").Delete()

# ------------------------------------------------------------------
# 5. Misc view bookkeeping to mirror the authored workbook.
# ------------------------------------------------------------------
$ws.Range("G6").Select()
